$d = $word.ActiveDocument
$d.Content.Find.Execute("February", $true, $false, $false, $false, $false,
                         $true, 1, $false, "March", 2)
